$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$lo  = $ws1.ListObjects.Item(1)

# --- fundamental_data: add a "currency" column right after ghg_s3 -----------
# Insert a blank column at M (13), pushing ghg_s3's right-hand neighbours
# (company_revenue .. company_cash_equivalents) one column to the right.
$ws1.Columns.Item(13).Insert()

# Populate the header and the two data rows for the new column.
$ws1.Range("M1").Value = "currency"
$ws1.Range("M2").Value = "USD"
$ws1.Range("M3").Value = "USD"

# The shifted header cells (now N1:R1) already carry the right text, but the
# table definition needs to be nudged to resync its column metadata with the
# actual header text, so re-assign each cell to itself.
foreach ($addr in @("N1", "O1", "P1", "Q1", "R1")) {
    $ws1.Range($addr).Value = $ws1.Range($addr).Value2
}

# Grow the table so the new column (and the now-18-column layout) is covered.
$lo.Resize($ws1.Range("A1:R51"))

# Re-sync one more time now that the table has actually grown to 18 columns,
# so the newly appended last column also reflects the true header text
# instead of a generic default name.
foreach ($addr in @("N1", "O1", "P1", "Q1", "R1")) {
    $ws1.Range($addr).Value = $ws1.Range($addr).Value2
}

# --- small (8pt) font used for the sheet's phonetic properties -------------
# Touch a scratch cell far outside the used range to register the new font
# in the style table, then remove every trace of the scratch cell itself so
# no visible cell ends up using it and the sheet dimension is unaffected.
$scratch = $ws1.Range("Z100")
$scratch.Font.Size = 8
$scratch.ClearFormats()
$scratch.ClearContents()
$scratch.EntireRow.Delete()

# --- selection / active sheet -----------------------------------------------
# fundamental_data becomes the active sheet with M4 selected; target_data
# keeps its own selection but is no longer the active tab.
$ws1.Activate()
$ws1.Range("M4").Select()
